$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '40.906.53'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.85%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.411.90'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.65%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '88.65'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.535'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.09%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.496'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.09%  '
$ws.Range('E10').Value = '  -2.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '31.80'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.07%  '
$ws.Range('E12').Value = '  -1.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.784.15'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.60%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.68'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.65'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.413.20'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.67%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.772'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.57%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '40.858.98'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0921'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.24'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.80'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.92'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '233.76'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.67'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.64%  '
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('E26').Value = '  -4.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.06%  '
$ws.Range('E28').Value = '  -2.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.52'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.25'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.25'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.92%  '
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('E33').Value = '  -5.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0741'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.44'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.89'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.87%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '16.52'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.09%  '
$ws.Range('E38').Value = '  -1.96%  '
$ws.Range('E39').Value = '  -5.40%  '
$ws.Range('E40').Value = '  -3.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.86'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.33'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.991.04'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '18.71'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.06%  '
$ws.Range('E45').Value = '  -4.60%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.87'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.36'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.650.52'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '94.00'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.25'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.56'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.66%  '
